# plots.xlsx — insert "Time(ms)" / "Size(mb)" header rows under each MODE
# label, shifting the existing data rows down (one new row per block).
#
# Layout before:
#   A1 = MODE1
#   A2:C11  = data block 1 (10 rows)
#   A13 = MODE2
#   A14:C23 = data block 2 (10 rows)
#
# Layout after:
#   A1 = MODE1
#   A2 = Size(mb) / C2 = Time(ms)      <- new header row
#   A3:C12  = data block 1 (shifted down by 1)
#   A14 = MODE2
#   A15 = Size(mb) / C15 = Time(ms)    <- new header row
#   A16:C25 = data block 2 (shifted down by 2 overall)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row right below each MODE label, pushing the data (and the
# second block's label) down. Doing the first insert shifts row 13 (MODE2)
# down to row 14, so the second insert target is row 15.
$ws.Rows("2:2").Insert()
$ws.Rows("15:15").Insert()

# Fill the two new header rows. Write the C-column text first so that
# "Time(ms)" lands before "Size(mb)" in the shared-strings table (matching
# the target uniqueCount/index order).
$ws.Range("C2").Value = "Time(ms)"
$ws.Range("A2").Value = "Size(mb)"
$ws.Range("C15").Value = "Time(ms)"
$ws.Range("A15").Value = "Size(mb)"

# The embedded chart's anchor needs to grow by one row's worth of height
# (its top edge moves down 1 row, its bottom edge moves down 2 rows) to
# stay aligned with the data block that now starts one row lower.
$rowHeight = $ws.Rows(1).RowHeight
$co = $ws.ChartObjects(1)
$co.Top = $co.Top + $rowHeight
$co.Height = $co.Height + $rowHeight

# Restore the active selection to match the saved workbook state.
$ws.Range("F26").Select()
